# Documentation/State.docx update:
# - Replace the "11/14/2016 Created a new class..." paragraph's text with
#   the new "You have to create a particle effect..." TODO entry (keeping
#   the existing _GoBack bookmark in that paragraph).
# - Insert a new "11/25/2016 Added the seal..." paragraph right after it.
# - Re-insert the original "11/14/2016 Created a new class..." paragraph
#   right after that one.

$d = $word.ActiveDocument

$target = $d.Paragraphs.Item(4)
$targetRange = $target.Range

# Create two new empty paragraphs immediately after the target paragraph.
# Calling InsertParagraphAfter twice on the (fixed) target's Range inserts
# both new paragraphs right after the target, in order.
$targetRange.InsertParagraphAfter()
$targetRange.InsertParagraphAfter()

# Fill in the two newly created paragraphs (now items 5 and 6).
$d.Paragraphs.Item(5).Range.Text = "11/25/2016 Added the seal to the left hand to teleport. The seal has the pointer invisible to know when it touches something. "
$d.Paragraphs.Item(6).Range.Text = "11/14/2016 Created a new class for Honovi that inherits from the abstract class Character. The came is still working."

# Finally, overwrite the original paragraph's text with the new TODO entry.
$targetRange.Text = "You have to create a particle effect for the platform that is the target of the teleportation. Make the particle effect disappear if the pointer leaves the platform. Make the seal active if you press the grip. Create a mask for the vision when teleporting."
